$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update data values in columns A and B (rows 1-32)
$ws.Range("A1").Value = -0.036339451647734222
$ws.Range("B1").Value = 0.036264386098480372
$ws.Range("A2").Value = 0.0098432878224379294
$ws.Range("B2").Value = -0.010134456652020063
$ws.Range("A3").Value = 0.11306239171234367
$ws.Range("B3").Value = -0.11340311284752858
$ws.Range("A4").Value = -0.17058276328828725
$ws.Range("B4").Value = 0.16985824376000735
$ws.Range("A5").Value = -0.16385824455573239
$ws.Range("B5").Value = 0.16240920654651969
$ws.Range("A6").Value = -0.054251703151948494
$ws.Range("B6").Value = 0.054216118310356443
$ws.Range("A7").Value = -0.034216119271624379
$ws.Range("B7").Value = 0.034163188511438847
$ws.Range("A8").Value = -0.014163189476651183
$ws.Range("B8").Value = 0.014135748232118495
$ws.Range("A9").Value = -0.0081357490623386042
$ws.Range("B9").Value = 0.0081126562404669755
$ws.Range("A10").Value = -0.0021126570722103111
$ws.Range("B10").Value = 0.0021135555315723309
$ws.Range("A11").Value = 0.0023864436513818532
$ws.Range("B11").Value = -0.0023952748906381771
$ws.Range("A12").Value = 0.0083952740589361419
$ws.Range("B12").Value = -0.0084627410758044697
$ws.Range("A13").Value = 0.014462740246029782
$ws.Range("B13").Value = -0.01449849453285168
$ws.Range("A14").Value = -0.032214379883003019
$ws.Range("B14").Value = 0.032175390571645757
$ws.Range("A15").Value = -0.02617539140321945
$ws.Range("B15").Value = 0.026141961209670939
$ws.Range("A16").Value = -0.020141962044320838
$ws.Range("B16").Value = 0.02010308898518165
$ws.Range("A17").Value = -0.014103089824286208
$ws.Range("B17").Value = 0.014090005254781879
$ws.Range("A18").Value = -0.036107674751360008
$ws.Range("B18").Value = 0.036095886414372558
$ws.Range("A19").Value = -0.027095887235272897
$ws.Range("B19").Value = 0.027013017870027856
$ws.Range("A20").Value = -0.018013018697848437
$ws.Range("B20").Value = 0.01800418390370595
$ws.Range("A21").Value = -0.0090041847325323943
$ws.Range("B21").Value = 0.0089999991704692306
$ws.Range("A22").Value = -0.093927935881005453
$ws.Range("B22").Value = 0.093621879404166108
$ws.Range("A23").Value = -0.08462188023135564
$ws.Range("B23").Value = 0.08412412734192376
$ws.Range("A24").Value = -0.042124128506798364
$ws.Range("B24").Value = 0.041999998829107099
$ws.Range("A25").Value = -0.10215750421659564
$ws.Range("B25").Value = 0.10201886157765472
$ws.Range("A26").Value = -0.096018862402864613
$ws.Range("B26").Value = 0.095843754749360244
$ws.Range("A27").Value = -0.089843755578891127
$ws.Range("B27").Value = 0.089258111012885166
$ws.Range("A28").Value = -0.083258111861181483
$ws.Range("B28").Value = 0.082872788394078079
$ws.Range("A29").Value = -0.070872789314401885
$ws.Range("B29").Value = 0.07070102473515405
$ws.Range("A30").Value = -0.050701025740301997
$ws.Range("B30").Value = 0.050385141211459761
$ws.Range("A31").Value = -0.035385142177137396
$ws.Range("B31").Value = 0.035297993251273141
$ws.Range("A32").Value = -0.014297994277663761
$ws.Range("B32").Value = 0.014197333089365216

# Widen column B to match column A (15.42578125 chars raw width).
# The engine quantizes ColumnWidth writes to a pixel grid, so feed it
# the character-width value whose rounded pixel width lands on the
# closest achievable grid point (15.5) to the target raw width.
$ws.Columns.Item(2).ColumnWidth = 14.665
